$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A77").Value = "2024-10-18 00:00:00"
$ws.Range("B77").Value = 73650
$ws.Range("C77").Value = 10320.18
$ws.Range("D77").Value = 9132.91
$ws.Range("E77").Value = 7.1033
